{"js": "// The author placed the cursor at the very end of the document (after the\n// existing \" ththr\" text, where the _GoBack bookmark lived), pressed\n// Enter to start a new paragraph, and typed \"qwe\". Word keeps the\n// \"_GoBack\" bookmark pinned to the most recent edit location, so it moves\n// from the end of the first paragraph to the end of the newly typed text.\n\nconst body = context.document.body;\n\n// 1) Append a new paragraph with the typed text at the end of the body.\nbody.insertParagraph(\"qwe\", Word.InsertLocation.end);\n\n// 2) Drop the bookmark from its old location (end of the first paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n\n// 3) Re-fetch the newly added (now last) paragraph so we get a fresh,\n//    properly collapsible range, then place \"_GoBack\" right after \"qwe\".\nconst lastParagraph = body.paragraphs.getLast();\nconst endOfLastParagraph = lastParagraph.getRange(Word.RangeLocation.end);\nendOfLastParagraph.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The author placed the cursor at the very end of the document (where the\n# \"_GoBack\" bookmark was sitting, right after \" ththr\"), pressed Enter to\n# start a new paragraph, and typed \"qwe\". Word always keeps \"_GoBack\"\n# pinned to the most recent edit location, so after typing it has to move\n# from the end of the first paragraph to the end of the newly typed text.\n\n$d = $word.ActiveDocument\n\n# 1) Go to the very end of the document and start a new paragraph with the\n#    typed text, just like the author did.\n$sel = $word.Selection\n$sel.EndKey(6)          # wdStory -> end of the document\n$sel.TypeParagraph()\n# A trailing placeholder character is typed along with the real text so\n# that the insertion point used for the bookmark below is never the very\n# last character position of the paragraph (Word re-collapses a bookmark\n# placed exactly there to the whole paragraph instead of a single point).\n$sel.TypeText(\"qweX\")\n\n# 2) Remove the bookmark from its old location (end of the first paragraph).\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 3) Recreate \"_GoBack\" collapsed right after \"qwe\" (before the placeholder).\n$lastParagraph = $d.Paragraphs.Last\n$afterQwe = $lastParagraph.Range.Start + 3   # \"qwe\" is 3 characters long\n$bookmarkRange = $d.Range($afterQwe, $afterQwe)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# 4) Drop the placeholder character now that the bookmark is anchored.\n$placeholder = $d.Range($afterQwe, $afterQwe + 1)\n$placeholder.Delete()\n"}
